$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings stored as text (e.g. "0.2960",
# "0.000007752", "31.123.82"). Force Text format first so Excel does not silently
# coerce these into numbers (which would drop trailing zeros / use sci notation).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "31.123.82"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").Value = "1.954.09"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "245.95"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "0.4895"
$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("D8").Value = "44.75"
$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("D9").Value = "0.2960"
$ws.Range("E9").Value = "  +1.41%  "

$ws.Range("D10").Value = "0.06830"
$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("D11").Value = "19.22"
$ws.Range("E11").Value = "  -1.18%  "

$ws.Range("D12").Value = "106.65"
$ws.Range("E12").Value = "  -5.83%  "

$ws.Range("D13").Value = "0.07731"
$ws.Range("E13").Value = "  +1.72%  "

$ws.Range("D14").Value = "1.928.74"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "5.411"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").Value = "0.7135"
$ws.Range("E16").Value = "  +4.83%  "

$ws.Range("D17").Value = "285.52"
$ws.Range("E17").Value = "  -4.41%  "

$ws.Range("D18").Value = "30.998.27"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("D19").Value = "0.000007752"
$ws.Range("E19").Value = "  +1.33%  "

$ws.Range("D20").Value = "13.21"

$ws.Range("B21").Value = "BitDAO"
$ws.Range("C21").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D21").Value = "0.4903"
$ws.Range("E21").Value = "  +9.56%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.180.53"
$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "5.517"
$ws.Range("E24").Value = "  -0.67%  "

$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "1.0000"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("B26").Value = "Chainlink"
$ws.Range("C26").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D26").Value = "6.605"
$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.921"
$ws.Range("E27").Value = "  +3.60%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "168.84"
$ws.Range("E28").Value = "  +0.13%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.98"
$ws.Range("E29").Value = "  -1.93%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "2.209"
$ws.Range("E30").Value = "  +3.88%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1054"
$ws.Range("E31").Value = "  -1.67%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "1.438"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.745"
$ws.Range("E33").Value = "  +15.58%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "4.477"
$ws.Range("E34").Value = "  +6.89%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.05006"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7624"
$ws.Range("E36").Value = "  +1.52%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.166"
$ws.Range("E37").Value = "  +1.46%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.729"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02044"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.705"
$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "2.149"
$ws.Range("E41").Value = "  +6.01%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.428"
$ws.Range("E42").Value = "  +10.24%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.4491"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "109.55"
$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.8805"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "72.63"
$ws.Range("E46").Value = "  +3.36%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9994"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.496"
$ws.Range("E48").Value = "  +2.14%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "981.47"
$ws.Range("E49").Value = "  +15.12%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1283"
$ws.Range("E50").Value = "  +3.94%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.302"
$ws.Range("E51").Value = "  -0.56%  "
